$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'332.03"
$ws.Range("E2").Value = "'0.67%"
$ws.Range("G2").Value = "'12"
$ws.Range("D3").Value = "'45.45"
$ws.Range("E3").Value = "'2.65%"
$ws.Range("G3").Value = "'12"
$ws.Range("D4").Value = "'5.493"
$ws.Range("E4").Value = "'-0.06%"
$ws.Range("G4").Value = "'12"
$ws.Range("D5").Value = "'0.08365"
$ws.Range("E5").Value = "'4.13%"
$ws.Range("G5").Value = "'12"
$ws.Range("D6").Value = "'2.092"
$ws.Range("E6").Value = "'1.56%"
$ws.Range("G6").Value = "'12"
$ws.Range("D7").Value = "'0.9889"
$ws.Range("E7").Value = "'3.50%"
$ws.Range("G7").Value = "'12"
$ws.Range("D8").Value = "'2.552"
$ws.Range("E8").Value = "'-2.89%"
$ws.Range("G8").Value = "'12"
$ws.Range("D9").Value = "'0.1200"
$ws.Range("E9").Value = "'4.98%"
$ws.Range("G9").Value = "'12"
$ws.Range("D10").Value = "'0.1935"
$ws.Range("E10").Value = "'2.19%"
$ws.Range("G10").Value = "'12"
$ws.Range("D11").Value = "'9.439"
$ws.Range("E11").Value = "'-6.44%"
$ws.Range("G11").Value = "'12"
$ws.Range("D12").Value = "'0.09833"
$ws.Range("E12").Value = "'-0.55%"
$ws.Range("G12").Value = "'12"
$ws.Range("D13").Value = "'0.04695"
$ws.Range("E13").Value = "'-3.97%"
$ws.Range("G13").Value = "'12"
$ws.Range("E14").Value = "'-0.52%"
$ws.Range("G14").Value = "'12"
$ws.Range("D15").Value = "'0.001286"
$ws.Range("E15").Value = "'2.20%"
$ws.Range("G15").Value = "'12"
$ws.Range("D16").Value = "'0.005954"
$ws.Range("E16").Value = "'-3.09%"
$ws.Range("G16").Value = "'12"
$ws.Range("D17").Value = "'3.392"
$ws.Range("E17").Value = "'0.10%"
$ws.Range("G17").Value = "'12"
$ws.Range("D18").Value = "'4.423"
$ws.Range("E18").Value = "'0.56%"
$ws.Range("G18").Value = "'12"
$ws.Range("E19").Value = "'-0.77%"
$ws.Range("G19").Value = "'12"
$ws.Range("D20").Value = "'0.1354"
$ws.Range("E20").Value = "'-1.99%"
$ws.Range("G20").Value = "'12"
$ws.Range("D21").Value = "'0.2543"
$ws.Range("E21").Value = "'-1.59%"
$ws.Range("G21").Value = "'12"
$ws.Range("D22").Value = "'0.04151"
$ws.Range("E22").Value = "'1.47%"
$ws.Range("G22").Value = "'12"
$ws.Range("D23").Value = "'0.001293"
$ws.Range("E23").Value = "'-0.40%"
$ws.Range("G23").Value = "'12"
$ws.Range("D24").Value = "'0.004574"
$ws.Range("E24").Value = "'4.94%"
$ws.Range("G24").Value = "'12"
$ws.Range("D25").Value = "'0.0001302"
$ws.Range("E25").Value = "'8.40%"
$ws.Range("G25").Value = "'12"
$ws.Range("E26").Value = "'-0.03%"
$ws.Range("G26").Value = "'12"
$ws.Range("G27").Value = "'12"
$ws.Range("G28").Value = "'12"
$ws.Range("G29").Value = "'12"
$ws.Range("G30").Value = "'12"
$ws.Range("G31").Value = "'12"
$ws.Range("G32").Value = "'12"
$ws.Range("G33").Value = "'12"
$ws.Range("G34").Value = "'12"
$ws.Range("G35").Value = "'12"
$ws.Range("G36").Value = "'12"
$ws.Range("G37").Value = "'12"
$ws.Range("D38").Value = "'0.02710"
$ws.Range("E38").Value = "'4.68%"
$ws.Range("G38").Value = "'12"
$ws.Range("D39").Value = "'0.05774"
$ws.Range("E39").Value = "'-0.36%"
$ws.Range("G39").Value = "'12"
$ws.Range("D40").Value = "'0.007893"
$ws.Range("E40").Value = "'4.27%"
$ws.Range("G40").Value = "'12"
$ws.Range("D41").Value = "'0.1433"
$ws.Range("E41").Value = "'1.99%"
$ws.Range("G41").Value = "'12"
$ws.Range("D42").Value = "'0.007816"
$ws.Range("E42").Value = "'6.73%"
$ws.Range("G42").Value = "'12"
$ws.Range("D43").Value = "'0.002102"
$ws.Range("E43").Value = "'4.26%"
$ws.Range("G43").Value = "'12"
$ws.Range("D44").Value = "'0.008950"
$ws.Range("E44").Value = "'-1.48%"
$ws.Range("G44").Value = "'12"
$ws.Range("D45").Value = "'0.3547"
$ws.Range("G45").Value = "'12"
$ws.Range("D46").Value = "'0.00007074"
$ws.Range("E46").Value = "'0.53%"
$ws.Range("G46").Value = "'12"
$ws.Range("E47").Value = "'0.08%"
$ws.Range("G47").Value = "'12"
$ws.Range("E48").Value = "'0.36%"
$ws.Range("G48").Value = "'12"
$ws.Range("D49").Value = "'0.003536"
$ws.Range("E49").Value = "'0.10%"
$ws.Range("G49").Value = "'12"
$ws.Range("D50").Value = "'0.003056"
$ws.Range("E50").Value = "'-13.37%"
$ws.Range("G50").Value = "'12"
$ws.Range("E51").Value = "'0.08%"
$ws.Range("G51").Value = "'12"
